$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Modify Transaction" sheet: just move the cached selection to L21.
# (Selecting it first, before the final Activate()/Select() below on
# "Transactions", keeps it out of the final tabSelected state.)
# ---------------------------------------------------------------------------
$wsModify = $wb.Worksheets.Item("Modify Transaction")
$wsModify.Range("L21").Select() | Out-Null

# ---------------------------------------------------------------------------
# "Transactions" sheet: entry-id / transaction-type / amount swap between
# rows 2 and 3 (a Waive Interest txn and a Repayment txn trade places and
# get new, higher entry ids), plus a new entry id on row 4.
# ---------------------------------------------------------------------------
$wsTxn = $wb.Worksheets.Item("Transactions")

# Row 2 becomes the "Repayment" transaction (10000 / 10000, #,##0 format).
$wsTxn.Range("A2").Value = 1052
$wsTxn.Range("D2").Value = "Repayment"
$wsTxn.Range("E2").Value = 10000
$wsTxn.Range("E2").NumberFormat = "#,##0"
$wsTxn.Range("F2").Value = 10000
$wsTxn.Range("F2").NumberFormat = "#,##0"
$wsTxn.Range("G2").Value = 0

# Row 3 becomes the "Waive interest" transaction (101.92 / general format).
$wsTxn.Range("A3").Value = 1051
$wsTxn.Range("D3").Value = "Waive interest"
$wsTxn.Range("E3").Value = 101.92
$wsTxn.Range("E3").ClearFormats()
$wsTxn.Range("E3").VerticalAlignment = -4108
$wsTxn.Range("F3").Value = 0
$wsTxn.Range("F3").ClearFormats()
$wsTxn.Range("F3").VerticalAlignment = -4108
$wsTxn.Range("G3").Value = 101.92
$wsTxn.Range("J3").Value = 10000
$wsTxn.Range("J3").NumberFormat = "#,##0"

# Row 4: only the entry id changes.
$wsTxn.Range("A4").Value = 1040

# Best-fit the Transaction Date / Transaction Type columns.
$wsTxn.Columns.Item(3).ColumnWidth = 15.140625
$wsTxn.Columns.Item(4).ColumnWidth = 15.28515625

# Make "Transactions" the active sheet/tab with C4 selected - this also
# clears tabSelected on whichever sheet previously had it ("Prepay Loan").
$wsTxn.Activate() | Out-Null
$wsTxn.Range("C4").Select() | Out-Null
